# Apply updated cryptocurrency price/volume data to Sheet1
# (commit: "Updated cryptos list on Sun Oct 27 19:24:26 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($sheet, $addr, $text) {
    $cell = $sheet.Range($addr)
    # Force text format so numeric-looking strings (e.g. "585.80", 
    # multi-dot "67.630.29", subscript-digit numbers) are kept verbatim
    # instead of being coerced into floating point numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

Set-TextCell $ws "D2" "67.630.29"
Set-TextCell $ws "E2" "  +0.78%  "
Set-TextCell $ws "D3" "2.487.01"
Set-TextCell $ws "E3" "  +0.19%  "
Set-TextCell $ws "E4" "  -0.01%  "
Set-TextCell $ws "D5" "585.80"
Set-TextCell $ws "E5" "  +0.08%  "
Set-TextCell $ws "D6" "176.41"
Set-TextCell $ws "E6" "  +2.53%  "
Set-TextCell $ws "E7" "  -0.02%  "
Set-TextCell $ws "E8" "  +0.19%  "
Set-TextCell $ws "E9" "  +3.35%  "
Set-TextCell $ws "E10" "  +0.10%  "
Set-TextCell $ws "E11" "  +2.56%  "
Set-TextCell $ws "E12" "  -0.04%  "
Set-TextCell $ws "D13" "2.940.37"
Set-TextCell $ws "E13" "  +0.21%  "
Set-TextCell $ws "D14" "25.61"
Set-TextCell $ws "E14" "  +0.68%  "
Set-TextCell $ws "D15" "67.486.14"
Set-TextCell $ws "E15" "  +0.68%  "
Set-TextCell $ws "E16" "  +0.38%  "
Set-TextCell $ws "D17" "2.485.86"
Set-TextCell $ws "E17" "  +0.63%  "
Set-TextCell $ws "D18" "7.52"
Set-TextCell $ws "E18" "  +1.67%  "
Set-TextCell $ws "D19" "10.92"
Set-TextCell $ws "E19" "  -0.69%  "
Set-TextCell $ws "D20" "349.94"
Set-TextCell $ws "E20" "  -0.09%  "
Set-TextCell $ws "E21" "  +2.07%  "
Set-TextCell $ws "E22" "  -0.22%  "
Set-TextCell $ws "D23" "70.62"
Set-TextCell $ws "E23" "  +3.06%  "
Set-TextCell $ws "E24" "  +1.11%  "
Set-TextCell $ws "E25" "  -2.65%  "
Set-TextCell $ws "D26" "9.09"
Set-TextCell $ws "E26" "  -1.77%  "
Set-TextCell $ws "D28" "0.998"
Set-TextCell $ws "E28" "  -0.16%  "
Set-TextCell $ws "D29" "0.0₃0898"
Set-TextCell $ws "E29" "  -0.15%  "
Set-TextCell $ws "D30" "504.94"
Set-TextCell $ws "E30" "  -0.74%  "
Set-TextCell $ws "D31" "7.78"
Set-TextCell $ws "E31" "  +1.20%  "
Set-TextCell $ws "E32" "  +2.17%  "
Set-TextCell $ws "E33" "  +0.23%  "
Set-TextCell $ws "E34" "  +0.00%  "
Set-TextCell $ws "E35" "  +4.28%  "
Set-TextCell $ws "D36" "162.22"
Set-TextCell $ws "E36" "  +1.85%  "
Set-TextCell $ws "E37" "  -0.12%  "
Set-TextCell $ws "D38" "18.28"
Set-TextCell $ws "E38" "  +0.29%  "
Set-TextCell $ws "E40" "  +0.05%  "
Set-TextCell $ws "E41" "  +3.63%  "
Set-TextCell $ws "D42" "0.328"
Set-TextCell $ws "E42" "  +0.31%  "
Set-TextCell $ws "E43" "  +0.27%  "
Set-TextCell $ws "E44" "  +0.97%  "
Set-TextCell $ws "D45" "144.57"
Set-TextCell $ws "E45" "  +1.18%  "
Set-TextCell $ws "E47" "  -0.19%  "
Set-TextCell $ws "B48" "BabyDogeCoin"
Set-TextCell $ws "C48" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextCell $ws "D48" "0.0₆0254"
Set-TextCell $ws "E48" "  +1.61%  "
Set-TextCell $ws "B49" "Cronos"
Set-TextCell $ws "C49" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell $ws "D49" "0.0742"
Set-TextCell $ws "E49" "  +1.86%  "
Set-TextCell $ws "B50" "Optimism"
Set-TextCell $ws "C50" "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
Set-TextCell $ws "D50" "1.58"
Set-TextCell $ws "E50" "  +0.81%  "
Set-TextCell $ws "B51" "Mantle"
Set-TextCell $ws "C51" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextCell $ws "D51" "0.585"
Set-TextCell $ws "E51" "  +0.32%  "
